$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 922.7143
$ws.Range("I6").Value = 87
$ws.Range("K6").Value = 261
$ws.Range("M6").Value = -149

# Row 9
$ws.Range("H9").Value = 151.5
$ws.Range("I9").Value = 153.66667
$ws.Range("K9").Value = 153.66667
$ws.Range("M9").Value = 15.33332999999999

# Row 21
$ws.Range("H21").Value = 11416.667
$ws.Range("I21").Value = 11416.667
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 11416.667
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -10948.667
$ws.Range("N21").ClearContents()

# Row 23
$ws.Range("H23").Value = 11416.667
$ws.Range("I23").Value = 11416.667
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 11416.667
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -11182.667
$ws.Range("N23").ClearContents()

# Row 43
$ws.Range("H43").Value = 3100
$ws.Range("I43").Value = 3133.3333
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 3133.3333
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -3064.3333
$ws.Range("N43").Value = -3138

# Row 88
$ws.Range("H88").Value = 2610.5
$ws.Range("J88").Value = 2461.4285
$ws.Range("L88").Value = 2461.4285
$ws.Range("N88").Value = -3273.4285

# Row 91
$ws.Range("H91").Value = 2610.5
$ws.Range("J91").Value = 2461.4285
$ws.Range("L91").Value = 2461.4285
$ws.Range("N91").Value = -5269.4285

# Row 100
$ws.Range("H100").Value = 4456.1
$ws.Range("I100").Value = 3426.8333
$ws.Range("J100").Value = 6000
$ws.Range("K100").Value = 3426.8333
$ws.Range("L100").Value = 6000
$ws.Range("M100").Value = -2885.8333
$ws.Range("N100").Value = -7082

# Row 111
$ws.Range("H111").Value = 555.8
$ws.Range("I111").Value = 555.8
$ws.Range("K111").Value = 1667.4
$ws.Range("M111").Value = 1399.6

# Row 113
$ws.Range("H113").Value = 3199.6
$ws.Range("I113").Value = 2749.5
$ws.Range("K113").Value = 2749.5
$ws.Range("M113").Value = 504.5

# Row 116
$ws.Range("H116").Value = 6238.5713
$ws.Range("I116").Value = 12187
$ws.Range("K116").Value = 12187
$ws.Range("M116").Value = -8745

# Row 118
$ws.Range("H118").Value = 965.3333
$ws.Range("I118").Value = 509.33334
$ws.Range("K118").Value = 1528.00002
$ws.Range("M118").Value = 128.9999800000001

# Row 138
$ws.Range("H138").Value = 3105.5
$ws.Range("I138").Value = 2099.1428
$ws.Range("J138").Value = 3476.2632
$ws.Range("K138").Value = 6297.428400000001
$ws.Range("L138").Value = 10428.7896
$ws.Range("M138").Value = -1157.428400000001
$ws.Range("N138").Value = -20708.7896

# Row 141
$ws.Range("H141").Value = 2835
$ws.Range("I141").Value = 2835
$ws.Range("K141").Value = 8505
$ws.Range("M141").Value = -3325

$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 3120.2666
$ws.Range("I132").Value = 3128.9285
$ws.Range("K132").Value = 9386.7855
$ws.Range("M132").Value = -6856.7855

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 376.13333
$ws.Range("I94").Value = 310.9
$ws.Range("J94").Value = 506.6
$ws.Range("K94").Value = 310.9
$ws.Range("L94").Value = 506.6
$ws.Range("M94").Value = 140.1
$ws.Range("N94").Value = -1408.6

# Row 134
$ws.Range("H134").Value = 3933.3333
$ws.Range("I134").Value = 3933.3333
$ws.Range("K134").Value = 11799.9999
$ws.Range("M134").Value = -9264.999899999999

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1693.8
$ws.Range("J22").Value = 3441.6667
$ws.Range("L22").Value = 3441.6667
$ws.Range("N22").Value = -4141.6667

# Row 37
$ws.Range("H37").Value = 90
$ws.Range("I37").Value = 90
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 90
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 17
$ws.Range("N37").ClearContents()

# Row 86
$ws.Range("H86").Value = 5333.3335
$ws.Range("I86").Value = 5000
$ws.Range("K86").Value = 5000
$ws.Range("M86").Value = -3877

# Row 89
$ws.Range("H89").Value = 5333.3335
$ws.Range("I89").Value = 5000
$ws.Range("K89").Value = 25000
$ws.Range("M89").Value = -19384

# Row 131
$ws.Range("H131").Value = 51374
$ws.Range("J131").Value = 51374
$ws.Range("L131").Value = 51374
$ws.Range("N131").Value = -61454

# Row 134
$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -12465

$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 100
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 300
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -127
$ws.Range("N16").ClearContents()

# Row 80
$ws.Range("H80").Value = 4695.4644
$ws.Range("I80").Value = 4729.476
$ws.Range("J80").Value = 4593.4287
$ws.Range("K80").Value = 14188.428
$ws.Range("L80").Value = 13780.2861
$ws.Range("M80").Value = -13252.428
$ws.Range("N80").Value = -15652.2861

# Row 83
$ws.Range("H83").Value = 4695.4644
$ws.Range("I83").Value = 4729.476
$ws.Range("J83").Value = 4593.4287
$ws.Range("K83").Value = 42565.284
$ws.Range("L83").Value = 41340.85830000001
$ws.Range("M83").Value = -37885.284
$ws.Range("N83").Value = -50700.85830000001

# Row 92
$ws.Range("H92").Value = 10001.5
$ws.Range("I92").Value = 10000
$ws.Range("J92").Value = 10003
$ws.Range("K92").Value = 30000
$ws.Range("L92").Value = 30009
$ws.Range("M92").Value = -28752
$ws.Range("N92").Value = -32505

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 4300
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940

# Row 131
$ws.Range("H131").Value = 40000
$ws.Range("I131").Value = 40000
$ws.Range("K131").Value = 40000
$ws.Range("M131").Value = -34960

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 4155.1055
$ws.Range("I46").Value = 3591.8333
$ws.Range("J46").Value = 4415.077
$ws.Range("K46").Value = 3591.8333
$ws.Range("L46").Value = 4415.077
$ws.Range("M46").Value = -3403.8333
$ws.Range("N46").Value = -4791.077

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1054.8667
$ws.Range("I96").Value = 981.36365
$ws.Range("J96").Value = 1257
$ws.Range("K96").Value = 981.36365
$ws.Range("L96").Value = 1257
$ws.Range("M96").Value = 391.63635
$ws.Range("N96").Value = -4003
